$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$row = 19

# Copy just the date cell's formatting from the row above so the new row's
# date reuses the existing date style instead of creating a new style entry.
$ws.Range("A18").Copy()
$ws.Range("A19").PasteSpecial(-4122) # xlPasteFormats

$ws.Cells.Item($row, 1).Value = 42622.885833333334
$ws.Cells.Item($row, 2).Value = 26
$ws.Cells.Item($row, 3).Value = 64
$ws.Cells.Item($row, 4).Value = 33
$ws.Cells.Item($row, 5).Value = 64
$ws.Cells.Item($row, 6).Value = 37
$ws.Cells.Item($row, 7).Value = 21008
$ws.Cells.Item($row, 8).Value = 18381
$ws.Cells.Item($row, 9).Value = 999
$ws.Cells.Item($row, 10).Value = 214
$ws.Cells.Item($row, 11).Value = 111
$ws.Cells.Item($row, 12).Value = 10
$ws.Cells.Item($row, 13).Value = 6
$ws.Cells.Item($row, 14).Value = "Named"
